$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Enable distinct Even-page and First-page headers/footers for this
# section. Touching the six HeaderFooter objects (Primary/Even/First for
# both Headers and Footers) mints the header2/3.xml + footer2/3.xml parts
# and wires up the six headerReference/footerReference entries in sectPr
# without flipping on DifferentFirstPageHeaderFooter / titlePg, matching
# the target layout. We must actually write into the even/first slots so
# the parts get created (merely probing .Exists does not mint them).

$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2
$wdHeaderFooterEvenPages = 3

$evenHeader = $sec.Headers($wdHeaderFooterEvenPages)
$evenHeader.Range.Text = ""

$evenFooter = $sec.Footers($wdHeaderFooterEvenPages)
$evenFooter.Range.Text = ""

$firstHeader = $sec.Headers($wdHeaderFooterFirstPage)
$firstHeader.Range.Text = ""

$firstFooter = $sec.Footers($wdHeaderFooterFirstPage)
$firstFooter.Range.Text = ""

# --- Bump the version string shown in the (now "default"/primary) header
# from "QAPYTH3 v3.1" to "QAPYTH3 v4".
$primaryHeader = $sec.Headers($wdHeaderFooterPrimary)
$primaryHeader.Range.Find.Execute("QAPYTH3 v3.1", $false, $false, $false, $false, $false, $true, 1, $false, "QAPYTH3 v4", 2)
